$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 930.3570999999999
$ws.Range("I92").Value = 1026.8572
$ws.Range("J92").Value = 833.8570999999999
$ws.Range("K92").Value = 1026.8572
$ws.Range("L92").Value = 833.8570999999999
$ws.Range("M92").Value = 221.1428000000001
$ws.Range("N92").Value = -3329.8571
$ws.Range("H107").Value = 703.5833
$ws.Range("I107").Value = 770.46155
$ws.Range("J107").Value = 529.7
$ws.Range("K107").Value = 770.46155
$ws.Range("L107").Value = 529.7
$ws.Range("M107").Value = 1149.53845
$ws.Range("N107").Value = -4369.7
$ws.Range("H113").Value = 2178.5715
$ws.Range("I113").Value = 2250
$ws.Range("J113").Value = 2150
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 2150
$ws.Range("M113").Value = 1004
$ws.Range("N113").Value = -8658
$ws.Range("H129").Value = 800.1539
$ws.Range("I129").Value = 503.5
$ws.Range("J129").Value = 932
$ws.Range("K129").Value = 1510.5
$ws.Range("L129").Value = 2796
$ws.Range("M129").Value = 3489.5
$ws.Range("N129").Value = -12796
$ws.Range("H137").Value = 17243298
$ws.Range("I137").Value = 1141.7894
$ws.Range("J137").Value = 50003396
$ws.Range("K137").Value = 3425.3682
$ws.Range("L137").Value = 150010188
$ws.Range("M137").Value = -875.3681999999999
$ws.Range("N137").Value = -150015288
$ws.Range("H138").Value = 2756.5454
$ws.Range("I138").Value = 2189.068
$ws.Range("J138").Value = 3513.182
$ws.Range("K138").Value = 6567.204000000001
$ws.Range("L138").Value = 10539.546
$ws.Range("M138").Value = -1427.204000000001
$ws.Range("N138").Value = -20819.546

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 22000
$ws.Range("J23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("N23").Value = -22518
$ws.Range("H32").Value = 2225.15
$ws.Range("I32").Value = 2225.15
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2225.15
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1938.15
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 3951.5715
$ws.Range("I61").Value = 3465.1428
$ws.Range("J61").Value = 4438
$ws.Range("K61").Value = 3465.1428
$ws.Range("L61").Value = 4438
$ws.Range("M61").Value = -3253.1428
$ws.Range("N61").Value = -4862
$ws.Range("H74").Value = 13370.833
$ws.Range("I74").Value = 1391.6666
$ws.Range("J74").Value = 25350
$ws.Range("K74").Value = 1391.6666
$ws.Range("L74").Value = 25350
$ws.Range("M74").Value = -517.6666
$ws.Range("N74").Value = -27098
$ws.Range("H77").Value = 13370.833
$ws.Range("I77").Value = 1391.6666
$ws.Range("J77").Value = 25350
$ws.Range("K77").Value = 6958.333000000001
$ws.Range("L77").Value = 126750
$ws.Range("M77").Value = -2590.333000000001
$ws.Range("N77").Value = -135486
$ws.Range("H97").Value = 626.1622
$ws.Range("I97").Value = 534.7742
$ws.Range("K97").Value = 534.7742
$ws.Range("M97").Value = -38.77419999999995
$ws.Range("H102").Value = 2287.375
$ws.Range("I102").Value = 2185.5715
$ws.Range("K102").Value = 2185.5715
$ws.Range("M102").Value = -563.5715
$ws.Range("H132").Value = 15715484
$ws.Range("I132").Value = 24553524
$ws.Range("K132").Value = 73660572
$ws.Range("M132").Value = -73658042
$ws.Range("H136").Value = 3951.5715
$ws.Range("I136").Value = 3465.1428
$ws.Range("J136").Value = 4438
$ws.Range("K136").Value = 10395.4284
$ws.Range("L136").Value = 13314
$ws.Range("M136").Value = -7845.428400000001
$ws.Range("N136").Value = -18414

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 247.07143
$ws.Range("I94").Value = 222.63637
$ws.Range("K94").Value = 222.63637
$ws.Range("M94").Value = 228.36363
$ws.Range("H107").Value = 886501
$ws.Range("I107").Value = 1285928.8
$ws.Range("J107").Value = 7760
$ws.Range("K107").Value = 1285928.8
$ws.Range("L107").Value = 7760
$ws.Range("M107").Value = -1284008.8
$ws.Range("N107").Value = -11600

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7579255.5
$ws.Range("I31").Value = 1524.909
$ws.Range("J31").Value = 15156986
$ws.Range("K31").Value = 1524.909
$ws.Range("L31").Value = 15156986
$ws.Range("M31").Value = -1229.909
$ws.Range("N31").Value = -15157576
$ws.Range("H34").Value = 7579255.5
$ws.Range("I34").Value = 1524.909
$ws.Range("J34").Value = 15156986
$ws.Range("K34").Value = 1524.909
$ws.Range("L34").Value = 15156986
$ws.Range("M34").Value = -1322.909
$ws.Range("N34").Value = -15157390
$ws.Range("H105").Value = 782.3333
$ws.Range("J105").Value = 1830
$ws.Range("L105").Value = 1830
$ws.Range("N105").Value = -5324

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1011.125
$ws.Range("I14").Value = 1011.125
$ws.Range("K14").Value = 3033.375
$ws.Range("M14").Value = -2860.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2082
$ws.Range("I97").Value = 1915
$ws.Range("J97").Value = 2750
$ws.Range("K97").Value = 1915
$ws.Range("L97").Value = 2750
$ws.Range("M97").Value = -1419
$ws.Range("N97").Value = -3742
$ws.Range("H122").Value = 5003.1
$ws.Range("I122").Value = 5810.4
$ws.Range("J122").Value = 2581.2
$ws.Range("K122").Value = 17431.2
$ws.Range("L122").Value = 7743.599999999999
$ws.Range("M122").Value = -14981.2
$ws.Range("N122").Value = -12643.6
$ws.Range("H132").Value = 62502904
$ws.Range("I132").Value = 111113240
$ws.Range("J132").Value = 3904
$ws.Range("K132").Value = 333339720
$ws.Range("L132").Value = 11712
$ws.Range("M132").Value = -333337190
$ws.Range("N132").Value = -16772

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 2187.2163
$ws.Range("I122").Value = 2042.4783
$ws.Range("J122").Value = 2425
$ws.Range("K122").Value = 6127.4349
$ws.Range("L122").Value = 7275
$ws.Range("M122").Value = -3677.4349
$ws.Range("N122").Value = -12175

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 296.125
$ws.Range("I113").Value = 287.85
$ws.Range("J113").Value = 337.5
$ws.Range("K113").Value = 863.5500000000001
$ws.Range("L113").Value = 1012.5
$ws.Range("M113").Value = 1306.45
$ws.Range("N113").Value = -5352.5
$ws.Range("H126").Value = 3709.1538
$ws.Range("I126").Value = 3221.9
$ws.Range("J126").Value = 5333.3335
$ws.Range("K126").Value = 9665.700000000001
$ws.Range("L126").Value = 16000.0005
$ws.Range("M126").Value = -7195.700000000001
$ws.Range("N126").Value = -20940.0005
$ws.Range("H132").Value = 35653150
$ws.Range("I132").Value = 49021708
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 147065124
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -147062594
$ws.Range("N132").Value = -16058.9999
